# E-gradanin sheet sortiraj pls
#
# This script normalises the Croatian-diacritic text in a few spots of the
# workbook and replaces the footnoted/rich-text year headers on the
# "Stanovnistvo" sheet with plain numeric years.
#
# Concretely:
#   1. Rename worksheet "Zivorođeni" -> "Zivorodeni" (drop diacritics) and
#      keep the Print_Titles defined name in sync with the new name.
#   2. Rename worksheet "E-građani" -> "E-gradani" (drop diacritics).
#   3. On "Stanovnistvo", replace the row-1 header cells B1:V1 (which held
#      text like "2001.1)" with a superscript footnote marker) with plain
#      numeric years 2001..2021.
#   4. On "Spol", fix cell C1 from "Zene" (with diacritic Z) to the
#      plain-ASCII "Zene" already used by D1.
#   5. Restore the on-screen selections to match where the editor left the
#      cursor on each touched sheet, without disturbing which tab is active.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Zivorođeni" -> "Zivorodeni" ---------------------------------
$wsZivorodeni = $wb.Worksheets.Item("Živorođeni")
$wsZivorodeni.Name = "Zivorodeni"

# Keep the Print_Titles defined name text in sync with the renamed sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Zivorodeni!Print_Titles") {
        $n.RefersTo = "=Zivorodeni!`$A:`$A,Zivorodeni!`$1:`$1"
    }
}

# --- 2. Rename "E-građani" -> "E-gradani" -----------------------------------
$wsEgradani = $wb.Worksheets.Item("E-građani")
$wsEgradani.Name = "E-gradani"

# --- 3. Stanovnistvo: plain numeric year headers ----------------------------
$wsStanovnistvo = $wb.Worksheets.Item("Stanovnistvo")
$years = 2001..2021
$col = 2
foreach ($y in $years) {
    $wsStanovnistvo.Cells.Item(1, $col).Value = $y
    $col = $col + 1
}

# --- 4. Spol: fix C1 text ----------------------------------------------------
$wsSpol = $wb.Worksheets.Item("Spol")
$wsSpol.Cells.Item(1, 3).Value = "Zene"

# --- 5. Restore selections ---------------------------------------------------
$wsStanovnistvo.Range("V1").Select()
$wsSpol.Range("F1").Select()
$wsEgradani.Range("K18").Select()
